$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 404, shifting the
# existing rows 404-424 down to 406-426 (dimension grows to A1:R426).
$ws.Rows("404:405").Insert()

# Row 404 - new weekly entry (Primera), Provincia del Elquí
$ws.Range("A404").Value = 7
$ws.Range("B404").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C404").Value = "Ñuble"
$ws.Range("D404").Value = 45147
$ws.Range("E404").Value = 16
$ws.Range("F404").Value = 100112017
$ws.Range("G404").Value = "Apio"
$ws.Range("H404").Value = "Americana (o)"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 180
$ws.Range("K404").Value = 6000
$ws.Range("L404").Value = 6000
$ws.Range("M404").Value = 6000
$ws.Range("N404").Value = "$/docena de matas"
$ws.Range("O404").Value = "Provincia del Elquí"
$ws.Range("P404").Value = 1000
$ws.Range("Q404").Value = 6
$ws.Range("R404").Value = "Hortaliza"

# Row 405 - new weekly entry (Segunda), Provincia del Elquí
$ws.Range("A405").Value = 7
$ws.Range("B405").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C405").Value = "Ñuble"
$ws.Range("D405").Value = 45147
$ws.Range("E405").Value = 16
$ws.Range("F405").Value = 100112017
$ws.Range("G405").Value = "Apio"
$ws.Range("H405").Value = "Americana (o)"
$ws.Range("I405").Value = "Segunda"
$ws.Range("J405").Value = 180
$ws.Range("K405").Value = 5000
$ws.Range("L405").Value = 5000
$ws.Range("M405").Value = 5000
$ws.Range("N405").Value = "$/docena de matas"
$ws.Range("O405").Value = "Provincia del Elquí"
$ws.Range("P405").Value = 833
$ws.Range("Q405").Value = 6
$ws.Range("R405").Value = "Hortaliza"
